# Regenerate handback report: new handoff GUIDs / xlf hashes / timestamps
# replace the previous batch's values across the Overview, zh-cn and de-de
# sheets (and their corresponding hyperlink display text).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---- new identifiers -------------------------------------------------
$file1      = "50301329-0412-4357-b78a-8319006306c1.md"
$file2      = "ffff8c8d31bd-be20-468e-bb80-044c1018d2c1.md"

$xlfZhCn    = "50301329-0412-4357-b78a-8319006306c1.c2814a41a5ddfb1a290a0e26663c54093882a4e0.zh-cn.xlf"
$xlfDeDe    = "50301329-0412-4357-b78a-8319006306c1.c2814a41a5ddfb1a290a0e26663c54093882a4e0.de-de.xlf"

$overviewDate = "2016-08-29 19:08:23"
$zhCnHoDate   = "2016-08-29 19:08:19"
$zhCnHbDate   = "2016-08-29 19:08:35"
$deDeHoDate   = "2016-08-29 19:08:23"
$deDeHbDate   = "2016-08-29 19:08:43"

# ---- Overview sheet ----------------------------------------------------
$ws1.Range("A2").Value2 = $file1
$ws1.Range("B2").Value2 = "e2e\" + $file1
$ws1.Range("G2").Value2 = $overviewDate

$ws1.Range("A3").Value2 = $file2
$ws1.Range("B3").Value2 = "e2e\" + $file2
$ws1.Range("G3").Value2 = $overviewDate

# ---- zh-cn sheet ---------------------------------------------------------
$ws2.Range("A2").Value2 = $file1
$ws2.Range("G2").Value2 = $xlfZhCn
$ws2.Range("H2").Value2 = $zhCnHoDate
$ws2.Range("I2").Value2 = $file1
$ws2.Range("J2").Value2 = $xlfZhCn
$ws2.Range("K2").Value2 = $zhCnHbDate

$ws2.Range("A3").Value2 = $file2
$ws2.Range("G3").Value2 = $xlfZhCn
$ws2.Range("H3").Value2 = $zhCnHoDate
$ws2.Range("I3").Value2 = $file2
$ws2.Range("J3").Value2 = $xlfZhCn
$ws2.Range("K3").Value2 = $zhCnHbDate

# ---- de-de sheet ---------------------------------------------------------
$ws3.Range("A2").Value2 = $file1
$ws3.Range("G2").Value2 = $xlfDeDe
$ws3.Range("H2").Value2 = $deDeHoDate
$ws3.Range("I2").Value2 = $file1
$ws3.Range("J2").Value2 = $xlfDeDe
$ws3.Range("K2").Value2 = $deDeHbDate

$ws3.Range("A3").Value2 = $file2
$ws3.Range("G3").Value2 = $xlfDeDe
$ws3.Range("H3").Value2 = $deDeHoDate
$ws3.Range("I3").Value2 = $file2
$ws3.Range("J3").Value2 = $xlfDeDe
$ws3.Range("K3").Value2 = $deDeHbDate

# ---- Hyperlinks ----------------------------------------------------------
# `Range.Hyperlinks.Delete()` clears every hyperlink on the sheet in this
# runtime, so rebuild each sheet's set in one pass (re-adding in the exact
# original order keeps the regenerated r:id sequence == the original one,
# since the target URLs/hosts themselves are unchanged by this edit).

$srcRepo = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b71b25a3cfcec084b8979e9f07bea0e2788a51d5/e2e/"
$zhcnRepo = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ae755bcbf54b1449f6aeb615ac139db6338c9dcb/e2e/"
$dedeRepo = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4521238c20c853fbaaa3fce3b8751eaaa5e73099/e2e/"

$origFile1 = "476c2da7-d610-499c-b355-48b333bf3e17.md"
$origFile2 = "54556c64-f107-42b1-9518-16d70dc87157.md"

# Overview: B2 -> rId2 (source repo, file1), B3 -> rId3 (source repo, file2)
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), ($srcRepo + $origFile1), "", "", ("e2e\" + $file1))
$ws1.Hyperlinks.Add($ws1.Range("B3"), ($srcRepo + $origFile2), "", "", ("e2e\" + $file2))

# zh-cn: A2 -> rId2 (source), I2 -> rId3 (zh-cn repo), A3 -> rId4 (source), I3 -> rId5 (zh-cn repo)
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), ($srcRepo + $origFile1), "", "", $file1)
$ws2.Hyperlinks.Add($ws2.Range("I2"), ($zhcnRepo + $origFile1), "", "", $file1)
$ws2.Hyperlinks.Add($ws2.Range("A3"), ($srcRepo + $origFile2), "", "", $file2)
$ws2.Hyperlinks.Add($ws2.Range("I3"), ($zhcnRepo + $origFile2), "", "", $file2)

# de-de: A2 -> rId2 (source), I2 -> rId3 (de-de repo), A3 -> rId4 (source), I3 -> rId5 (de-de repo)
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), ($srcRepo + $origFile1), "", "", $file1)
$ws3.Hyperlinks.Add($ws3.Range("I2"), ($dedeRepo + $origFile1), "", "", $file1)
$ws3.Hyperlinks.Add($ws3.Range("A3"), ($srcRepo + $origFile2), "", "", $file2)
$ws3.Hyperlinks.Add($ws3.Range("I3"), ($dedeRepo + $origFile2), "", "", $file2)
